$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1278.2559
$ws.Range("I15").Value = 1278.2559
$ws.Range("K15").Value = 3834.7677
$ws.Range("M15").Value = -3665.7677
$ws.Range("H17").Value = 1463.7046
$ws.Range("J17").Value = 1463.7046
$ws.Range("L17").Value = 4391.1138
$ws.Range("N17").Value = -4727.1138
$ws.Range("H40").Value = 17254668
$ws.Range("I40").Value = 15979.733
$ws.Range("J40").Value = 35724692
$ws.Range("K40").Value = 15979.733
$ws.Range("L40").Value = 35724692
$ws.Range("M40").Value = -15804.733
$ws.Range("N40").Value = -35725042
$ws.Range("H61").Value = 175
$ws.Range("J61").Value = 17
$ws.Range("L61").Value = 51
$ws.Range("N61").Value = -395
$ws.Range("H98").Value = 2457.6
$ws.Range("I98").Value = 1905.3334
$ws.Range("K98").Value = 1905.3334
$ws.Range("M98").Value = -407.3334
$ws.Range("H106").Value = 47621904
$ws.Range("I106").Value = 55557220
$ws.Range("K106").Value = 55557220
$ws.Range("M106").Value = -55556589
$ws.Range("H122").Value = 2457.6
$ws.Range("I122").Value = 1905.3334
$ws.Range("K122").Value = 5716.0002
$ws.Range("M122").Value = -3266.0002
$ws.Range("H132").Value = 12717.519
$ws.Range("I132").Value = 3779.2
$ws.Range("K132").Value = 11337.6
$ws.Range("M132").Value = -8807.599999999999
$ws.Range("H137").Value = 13209350
$ws.Range("I137").Value = 1000879.2
$ws.Range("K137").Value = 3002637.6
$ws.Range("M137").Value = -3000087.6
$ws.Range("H138").Value = 2578.75
$ws.Range("I138").Value = 1941.7646
$ws.Range("J138").Value = 2979.8147
$ws.Range("K138").Value = 5825.293799999999
$ws.Range("L138").Value = 8939.444100000001
$ws.Range("M138").Value = -685.2937999999995
$ws.Range("N138").Value = -19219.4441

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10340.011
$ws.Range("J32").Value = 9902.125
$ws.Range("L32").Value = 9902.125
$ws.Range("N32").Value = -10476.125
$ws.Range("H110").Value = 1707755.6
$ws.Range("I110").Value = 2554396
$ws.Range("K110").Value = 2554396
$ws.Range("M110").Value = -2552351

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6097.755
$ws.Range("I31").Value = 1490.6428
$ws.Range("J31").Value = 7940.6
$ws.Range("K31").Value = 1490.6428
$ws.Range("L31").Value = 7940.6
$ws.Range("M31").Value = -1195.6428
$ws.Range("N31").Value = -8530.6
$ws.Range("H34").Value = 6097.755
$ws.Range("I34").Value = 1490.6428
$ws.Range("J34").Value = 7940.6
$ws.Range("K34").Value = 1490.6428
$ws.Range("L34").Value = 7940.6
$ws.Range("M34").Value = -1288.6428
$ws.Range("N34").Value = -8344.6
$ws.Range("H58").Value = 281227.75
$ws.Range("I58").Value = 528121.5600000001
$ws.Range("K58").Value = 528121.5600000001
$ws.Range("M58").Value = -527918.5600000001
$ws.Range("H86").Value = 7296.3184
$ws.Range("J86").Value = 10589.6
$ws.Range("L86").Value = 10589.6
$ws.Range("N86").Value = -12835.6
$ws.Range("H89").Value = 7296.3184
$ws.Range("J89").Value = 10589.6
$ws.Range("L89").Value = 52948
$ws.Range("N89").Value = -64180
$ws.Range("H132").Value = 13906070
$ws.Range("I132").Value = 16681607
$ws.Range("J132").Value = 28384.75
$ws.Range("K132").Value = 50044821
$ws.Range("L132").Value = 85154.25
$ws.Range("M132").Value = -50042291
$ws.Range("N132").Value = -90214.25
$ws.Range("H136").Value = 281227.75
$ws.Range("I136").Value = 528121.5600000001
$ws.Range("K136").Value = 1584364.68
$ws.Range("M136").Value = -1581814.68
$ws.Range("H138").Value = 69342.69
$ws.Range("J138").Value = 64963.145
$ws.Range("L138").Value = 64963.145
$ws.Range("N138").Value = -75243.14499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 357.5
$ws.Range("I86").Value = 315
$ws.Range("K86").Value = 945
$ws.Range("M86").Value = 241
$ws.Range("H89").Value = 357.5
$ws.Range("I89").Value = 315
$ws.Range("K89").Value = 2835
$ws.Range("M89").Value = 3093
$ws.Range("H108").Value = 13747.444
$ws.Range("I108").Value = 1242.3334
$ws.Range("K108").Value = 3727.0002
$ws.Range("M108").Value = -847.0001999999999
$ws.Range("H109").Value = 10152.5
$ws.Range("I109").Value = 4305
$ws.Range("K109").Value = 12915
$ws.Range("M109").Value = -11875
$ws.Range("H115").Value = 10899.8
$ws.Range("I115").Value = 4999.5
$ws.Range("K115").Value = 14998.5
$ws.Range("M115").Value = -13823.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3794034.5
$ws.Range("I70").Value = 4551041.5
$ws.Range("K70").Value = 4551041.5
$ws.Range("M70").Value = -4550771.5
$ws.Range("H73").Value = 3794034.5
$ws.Range("I73").Value = 4551041.5
$ws.Range("K73").Value = 4551041.5
$ws.Range("M73").Value = -4550105.5
$ws.Range("H97").Value = 534.7895
$ws.Range("I97").Value = 434.2414
$ws.Range("K97").Value = 434.2414
$ws.Range("M97").Value = 61.7586
$ws.Range("H132").Value = 90926.125
$ws.Range("I132").Value = 112222.69
$ws.Range("J132").Value = 9999.200000000001
$ws.Range("K132").Value = 336668.07
$ws.Range("L132").Value = 29997.6
$ws.Range("M132").Value = -334138.07
$ws.Range("N132").Value = -35057.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2833.2632
$ws.Range("I93").Value = 2879.5557
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 2879.5557
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -1631.5557
$ws.Range("N93").Value = -4496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5498.3076
$ws.Range("I62").Value = 3684.875
$ws.Range("J62").Value = 8399.799999999999
$ws.Range("K62").Value = 3684.875
$ws.Range("L62").Value = 8399.799999999999
$ws.Range("M62").Value = -3060.875
$ws.Range("N62").Value = -9647.799999999999
$ws.Range("H65").Value = 5498.3076
$ws.Range("I65").Value = 3684.875
$ws.Range("J65").Value = 8399.799999999999
$ws.Range("K65").Value = 18424.375
$ws.Range("L65").Value = 41999
$ws.Range("M65").Value = -15304.375
$ws.Range("N65").Value = -48239
$ws.Range("H113").Value = 913.2
$ws.Range("I113").Value = 892.9
$ws.Range("K113").Value = 2678.7
$ws.Range("M113").Value = -508.6999999999998
$ws.Range("H136").Value = 7326.28
$ws.Range("I136").Value = 1836.7646
$ws.Range("K136").Value = 5510.293799999999
$ws.Range("M136").Value = -2960.293799999999
$ws.Range("H140").Value = 81627.336
$ws.Range("J140").Value = 82246
$ws.Range("L140").Value = 82246
$ws.Range("N140").Value = -92606
$ws.Range("H141").Value = 29999.5
$ws.Range("J141").Value = 29999.5
$ws.Range("L141").Value = 29999.5
$ws.Range("N141").Value = -40359.5
